$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 50000
$ws.Range("J3").Value = 50000
$ws.Range("L3").Value = 50000
$ws.Range("N3").Value = -50228
$ws.Range("H15").Value = 1439.6052
$ws.Range("I15").Value = 1439.6052
$ws.Range("K15").Value = 4318.8156
$ws.Range("M15").Value = -4149.8156
$ws.Range("H33").Value = 357.3
$ws.Range("I33").Value = 266.0625
$ws.Range("J33").Value = 722.25
$ws.Range("K33").Value = 266.0625
$ws.Range("L33").Value = 722.25
$ws.Range("M33").Value = -37.0625
$ws.Range("N33").Value = -1180.25
$ws.Range("H70").Value = 1221022.1
$ws.Range("J70").Value = 2017.1666
$ws.Range("L70").Value = 6051.4998
$ws.Range("N70").Value = -6591.4998
$ws.Range("H73").Value = 1221022.1
$ws.Range("J73").Value = 2017.1666
$ws.Range("L73").Value = 6051.4998
$ws.Range("N73").Value = -7923.4998
$ws.Range("H98").Value = 3290619.5
$ws.Range("I98").Value = 3677475.2
$ws.Range("K98").Value = 3677475.2
$ws.Range("M98").Value = -3675977.2
$ws.Range("H102").Value = 50000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -56490
$ws.Range("H103").Value = 62501310
$ws.Range("I103").Value = 750
$ws.Range("K103").Value = 2250
$ws.Range("M103").Value = -1664
$ws.Range("H122").Value = 3290619.5
$ws.Range("I122").Value = 3677475.2
$ws.Range("K122").Value = 11032425.6
$ws.Range("M122").Value = -11029975.6
$ws.Range("H132").Value = 5879.2583
$ws.Range("I132").Value = 2693.2
$ws.Range("K132").Value = 8079.599999999999
$ws.Range("M132").Value = -5549.599999999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5804.6
$ws.Range("J2").Value = 3675
$ws.Range("L2").Value = 3675
$ws.Range("N2").Value = -3901
$ws.Range("H32").Value = 3997.5122
$ws.Range("I32").Value = 4206.4243
$ws.Range("K32").Value = 4206.4243
$ws.Range("M32").Value = -3919.4243
$ws.Range("H61").Value = 9136450
$ws.Range("I61").Value = 11768439
$ws.Range("K61").Value = 11768439
$ws.Range("M61").Value = -11768227
$ws.Range("H74").Value = 1706.1034
$ws.Range("I74").Value = 1166.56
$ws.Range("K74").Value = 1166.56
$ws.Range("M74").Value = -292.5599999999999
$ws.Range("H77").Value = 1706.1034
$ws.Range("I77").Value = 1166.56
$ws.Range("K77").Value = 5832.799999999999
$ws.Range("M77").Value = -1464.799999999999
$ws.Range("H110").Value = 5650.1177
$ws.Range("I110").Value = 5218
$ws.Range("K110").Value = 5218
$ws.Range("M110").Value = -3173
$ws.Range("H116").Value = 5804.6
$ws.Range("J116").Value = 3675
$ws.Range("L116").Value = 3675
$ws.Range("N116").Value = -8263
$ws.Range("H122").Value = 3987.4285
$ws.Range("I122").Value = 4470.7144
$ws.Range("K122").Value = 13412.1432
$ws.Range("M122").Value = -10962.1432
$ws.Range("H132").Value = 1669479.5
$ws.Range("I132").Value = 2658.4285
$ws.Range("J132").Value = 25004974
$ws.Range("K132").Value = 7975.2855
$ws.Range("L132").Value = 75014922
$ws.Range("M132").Value = -5445.2855
$ws.Range("N132").Value = -75019982
$ws.Range("H136").Value = 9136450
$ws.Range("I136").Value = 11768439
$ws.Range("K136").Value = 35305317
$ws.Range("M136").Value = -35302767

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5804.6
$ws.Range("J3").Value = 3675
$ws.Range("L3").Value = 3675
$ws.Range("N3").Value = -3903
$ws.Range("H107").Value = 3936.3333
$ws.Range("J107").Value = 6999
$ws.Range("L107").Value = 6999
$ws.Range("N107").Value = -10839

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 11118458
$ws.Range("I16").Value = 33337332
$ws.Range("K16").Value = 33337332
$ws.Range("M16").Value = -33337045
$ws.Range("H31").Value = 30305414
$ws.Range("I31").Value = 41668948
$ws.Range("J31").Value = 2655.6667
$ws.Range("K31").Value = 41668948
$ws.Range("L31").Value = 2655.6667
$ws.Range("M31").Value = -41668653
$ws.Range("N31").Value = -3245.6667
$ws.Range("H34").Value = 30305414
$ws.Range("I34").Value = 41668948
$ws.Range("J34").Value = 2655.6667
$ws.Range("K34").Value = 41668948
$ws.Range("L34").Value = 2655.6667
$ws.Range("M34").Value = -41668746
$ws.Range("N34").Value = -3059.6667
$ws.Range("H113").Value = 11118458
$ws.Range("I113").Value = 33337332
$ws.Range("K113").Value = 33337332
$ws.Range("M113").Value = -33335162
$ws.Range("H132").Value = 2860.926
$ws.Range("I132").Value = 2761.7144
$ws.Range("K132").Value = 8285.143199999999
$ws.Range("M132").Value = -5755.143199999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5978511
$ws.Range("I4").Value = 5978511
$ws.Range("K4").Value = 17935533
$ws.Range("M4").Value = -17935421
$ws.Range("H130").Value = 12219.111
$ws.Range("J130").Value = 15222
$ws.Range("L130").Value = 45666
$ws.Range("N130").Value = -55706

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 39000
$ws.Range("H80").Value = 4935.1113
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 5202.2856
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 5202.2856
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -7198.2856
$ws.Range("H83").Value = 4935.1113
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 5202.2856
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 26011.428
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -35995.428
$ws.Range("H102").Value = 3084.1428
$ws.Range("I102").Value = 3115.1667
$ws.Range("K102").Value = 3115.1667
$ws.Range("M102").Value = -1493.1667
$ws.Range("H132").Value = 8335958.5
$ws.Range("I132").Value = 2863.9092
$ws.Range("J132").Value = 100000000
$ws.Range("K132").Value = 8591.7276
$ws.Range("L132").Value = 300000000
$ws.Range("M132").Value = -6061.7276
$ws.Range("N132").Value = -300005060

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6401.467
$ws.Range("I7").Value = 5107.091
$ws.Range("K7").Value = 5107.091
$ws.Range("M7").Value = -4995.091
$ws.Range("H61").Value = 62504740
$ws.Range("I61").Value = 125001070
$ws.Range("J61").Value = 8406.5
$ws.Range("K61").Value = 125001070
$ws.Range("L61").Value = 8406.5
$ws.Range("M61").Value = -125000868
$ws.Range("N61").Value = -8810.5
$ws.Range("H68").Value = 5210767
$ws.Range("I68").Value = 8335269
$ws.Range("J68").Value = 3264.3333
$ws.Range("K68").Value = 8335269
$ws.Range("L68").Value = 3264.3333
$ws.Range("M68").Value = -8334520
$ws.Range("N68").Value = -4762.3333
$ws.Range("H71").Value = 5210767
$ws.Range("I71").Value = 8335269
$ws.Range("J71").Value = 3264.3333
$ws.Range("K71").Value = 41676345
$ws.Range("L71").Value = 16321.6665
$ws.Range("M71").Value = -41672601
$ws.Range("N71").Value = -23809.6665
$ws.Range("H113").Value = 62504740
$ws.Range("I113").Value = 125001070
$ws.Range("J113").Value = 8406.5
$ws.Range("K113").Value = 125001070
$ws.Range("L113").Value = 8406.5
$ws.Range("M113").Value = -124998900
$ws.Range("N113").Value = -12746.5
$ws.Range("H126").Value = 6401.467
$ws.Range("I126").Value = 5107.091
$ws.Range("K126").Value = 15321.273
$ws.Range("M126").Value = -12851.273
$ws.Range("H132").Value = 3188.6667
$ws.Range("J132").Value = 7999.4
$ws.Range("L132").Value = 23998.2
$ws.Range("N132").Value = -29058.2

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 99999
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H52").Value = 28362.666
$ws.Range("I52").Value = 30042
$ws.Range("K52").Value = 30042
$ws.Range("M52").Value = -29816
$ws.Range("H96").Value = 14896
$ws.Range("J96").Value = 18099.666
$ws.Range("L96").Value = 18099.666
$ws.Range("N96").Value = -20845.666
$ws.Range("H113").Value = 687.6667
$ws.Range("I113").Value = 785
$ws.Range("K113").Value = 2355
$ws.Range("M113").Value = -185
$ws.Range("H121").Value = 98439.664
$ws.Range("I121").Value = 99900
$ws.Range("J121").Value = 97709.5
$ws.Range("K121").Value = 99900
$ws.Range("L121").Value = 97709.5
$ws.Range("M121").Value = -98153
$ws.Range("N121").Value = -101203.5
$ws.Range("H122").Value = 3139
$ws.Range("I122").Value = 2405.8333
$ws.Range("K122").Value = 7217.499899999999
$ws.Range("M122").Value = -4767.499899999999
$ws.Range("H126").Value = 6124.303
$ws.Range("I126").Value = 7093.5
$ws.Range("J126").Value = 2524.4285
$ws.Range("K126").Value = 21280.5
$ws.Range("L126").Value = 7573.2855
$ws.Range("M126").Value = -18810.5
$ws.Range("N126").Value = -12513.2855
$ws.Range("H132").Value = 359189.78
$ws.Range("I132").Value = 2038
$ws.Range("K132").Value = 6114
$ws.Range("M132").Value = -3584
$ws.Range("H136").Value = 305630.44
$ws.Range("I136").Value = 2446.6333
$ws.Range("J136").Value = 3337468.2
$ws.Range("K136").Value = 7339.8999
$ws.Range("L136").Value = 10012404.6
$ws.Range("M136").Value = -4789.8999
$ws.Range("N136").Value = -10017504.6
